$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.300.51"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "1.786.72"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'338.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.3828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.3457"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'47.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'1.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "'0.07444"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'23.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.56%  "
$ws.Range("D13").Value = "'0.9999"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'6.452"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'7.281"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "1.783.88"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'0.00001078"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "'0.06651"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'82.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'17.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'6.469"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "28.289.04"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "'2.372"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "'1.446"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.430"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'155.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "'137.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").Value = "1.986.54"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'6.187"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("D33").Value = "'3.940"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "'0.08915"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").Value = "'0.02444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("D37").Value = "'0.6881"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'0.06387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "'0.2182"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'1.243"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "'1.504"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.33%  "
$ws.Range("D43").Value = "'8.321"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "'14.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("D45").Value = "'0.9988"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'0.6325"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").Value = "'3.873"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'133.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "'2.098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("D50").Value = "'0.07479"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("D51").Value = "'1.226"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.95%  "
